# Applies the 2023-12-02 20:45 "Atualizado por script" update to the
# netherlands_eredivisie_2023-2024 sheet:
#   1) Eight pairs of adjacent rows (same kick-off date/time) have their
#      match data (columns F:V) swapped - the scraper re-ordered the two
#      fixtures played at the same time on the same matchday.
#   2) Four brand-new fixture rows (118-121, Indice 117-120) are appended
#      after the previous last row (117), extending the used range from
#      A1:V117 to A1:V121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

# --- 1) swap F:V between each pair of rows -------------------------------
$swapPairs = @(
    @(56, 57),
    @(69, 70),
    @(75, 76),
    @(87, 88),
    @(104, 105),
    @(109, 110),
    @(111, 112),
    @(114, 115)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $cell1 = $col + $r1
        $cell2 = $col + $r2
        $v1 = $ws.Range($cell1).Value2
        $v2 = $ws.Range($cell2).Value2
        $ws.Range($cell1).Value2 = $v2
        $ws.Range($cell2).Value2 = $v1
    }
}

# --- 2) append the four new fixture rows ---------------------------------
$newRows = @(
    @{ Row=118; A=117; E=45262.78125;        F="FC Volendam"; G=0; H="Zwolle";          I=5; J=2.63; K="26/11/2023 14:42"; L=2.82; M="02/12/2023 18:36"; N=3.72; O="26/11/2023 14:42"; P=3.85; Q="02/12/2023 18:36"; R=2.61; S="26/11/2023 14:42"; T=2.42; U="02/12/2023 18:36"; V="https://www.betexplorer.com/football/netherlands/eredivisie/fc-volendam-zwolle/Ecc7ZbYI/" },
    @{ Row=119; A=118; E=45262.83333333334;  F="Sittard";     G=3; H="Vitesse";         I=1; J=2.03; K="25/11/2023 21:12"; L=2.06; M="02/12/2023 19:34"; N=3.74; O="25/11/2023 21:12"; P=3.57; Q="02/12/2023 19:34"; R=3.66; S="25/11/2023 21:12"; T=3.79; U="02/12/2023 19:34"; V="https://www.betexplorer.com/football/netherlands/eredivisie/sittard-vitesse/jFtjbrBJ/" },
    @{ Row=120; A=119; E=45262.83333333334;  F="Waalwijk";    G=2; H="Excelsior";       I=2; J=1.93; K="25/11/2023 21:12"; L=2.04; M="02/12/2023 19:36"; N=4.01; O="25/11/2023 21:12"; P=3.78; Q="02/12/2023 19:58"; R=3.68; S="25/11/2023 21:12"; T=3.57; U="02/12/2023 19:36"; V="https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-excelsior/U5ufc2QP/" },
    @{ Row=121; A=120; E=45262.875;          F="Heracles";    G=0; H="Sparta Rotterdam"; I=1; J=2.59; K="26/11/2023 17:12"; L=2.64; M="02/12/2023 20:57"; N=3.54; O="26/11/2023 17:12"; P=3.62; Q="02/12/2023 20:57"; R=2.76; S="26/11/2023 17:12"; T=2.68; U="02/12/2023 20:46"; V="https://www.betexplorer.com/football/netherlands/eredivisie/heracles-sparta-rotterdam/z7b3zuJC/" }
)

$ws.Range("A2").Copy()
foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").PasteSpecial(-4122)
}

$ws.Range("E2").Copy()
foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("E$r").PasteSpecial(-4122)
}

foreach ($nr in $newRows) {
    $r = $nr.Row

    $ws.Cells.Item($r, 1).Value2 = $nr.A

    $ws.Cells.Item($r, 2).Value2 = "netherlands"
    $ws.Cells.Item($r, 3).Value2 = "eredivisie"
    $ws.Cells.Item($r, 4).Value2 = "2023-2024"

    $ws.Cells.Item($r, 5).Value2 = $nr.E

    $ws.Range("F$r").Value2 = $nr.F
    $ws.Range("G$r").Value2 = $nr.G
    $ws.Range("H$r").Value2 = $nr.H
    $ws.Range("I$r").Value2 = $nr.I
    $ws.Range("J$r").Value2 = $nr.J
    $ws.Range("K$r").Value2 = $nr.K
    $ws.Range("L$r").Value2 = $nr.L
    $ws.Range("M$r").Value2 = $nr.M
    $ws.Range("N$r").Value2 = $nr.N
    $ws.Range("O$r").Value2 = $nr.O
    $ws.Range("P$r").Value2 = $nr.P
    $ws.Range("Q$r").Value2 = $nr.Q
    $ws.Range("R$r").Value2 = $nr.R
    $ws.Range("S$r").Value2 = $nr.S
    $ws.Range("T$r").Value2 = $nr.T
    $ws.Range("U$r").Value2 = $nr.U
    $ws.Range("V$r").Value2 = $nr.V
}
